$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 82 is a blank separator row (same style as other separator rows, e.g. row 71)
$ws.Range("A71:C71").Copy() | Out-Null
$ws.Range("A82").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New timesheet entries for Feb 07 2020 (rows 83-91)
$data = @(
    @("Feb 07 10:00 to 11:00", "Client call", "Sapphire Automation"),
    @("Feb 07 11:00 to 12:00", "Practicing neural network on local machine.", "Infimetrics"),
    @("Feb 07 12:00 to 13:00", "Used OOP concepts of inheritence and polymorphism in MLP code.", "Infimetrics"),
    @("Feb 07 13:00 to 14:00", "Lunch", "Infimetrics"),
    @("Feb 07 14:00 to 15:00", "Modified logic of dot product in MLP.", "Infimetrics"),
    @("Feb 07 15:00 to 16:00", "Implementing back propagation in MLP.", "Infimetrics"),
    @("Feb 07 16:00 to 17:00", "Succesfully implemented MLP", "Infimetrics"),
    @("Feb 07 17:00 to 18:00", "Verified model accuracy", "Infimetrics"),
    @("Feb 07 18:00 to 19:00", "Model accuracy is very poor, need modification for algorithm.", "Infimetrics")
)

$startRow = 83
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowData = $data[$i]

    # Copy formatting/style from an existing "normal" entry row (row 72) that uses styles 1/3/1
    $ws.Range("A72").Copy() | Out-Null
    $ws.Range("A$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws.Range("B72").Copy() | Out-Null
    $ws.Range("B$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("C72").Copy() | Out-Null
    $ws.Range("C$r").PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
}

$excel.CutCopyMode = 0

# Excel recalculated the autofit row heights of a couple of wrapped-text rows
# once the new content/rows were added further down the sheet.
$ws.Rows.Item(73).RowHeight = 45
$ws.Rows.Item(77).RowHeight = 30

$ws.Application.ActiveWindow.ScrollRow = 82
$ws.Range("D91").Select()
